$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store these numeric-looking values as
# text (matching the original inlineStr/text cell type) instead of
# auto-converting them to numbers.
$ws.Range("A1").Value = "'1.5"
$ws.Range("B1").Value = "'-2.45"
$ws.Range("A2").Value = "'-1"
$ws.Range("B2").Value = "'3.8"
